# Standard User name change in Contacts test data files - 6th Mar 2024
# Update the StdUser value on the "Users" sheet from "Drew Koecher" to "Ayati Arvind".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Ayati Arvind"
